$wb = $excel.ActiveWorkbook
$bingo = $wb.Worksheets.Item("bingo_data")

# Re-enter G2:G33 as a single range formula so Excel groups it into one shared-formula
# block (matching the author re-filling the y-axis-label column down from G2).
$bingo.Range("G2:G33").Formula = "=(F2-`$I`$2)/(`$I`$3-`$I`$2)*90+10"

# Move bingo_data's own selection (it will no longer be the active sheet)
$bingo.Range("I9").Select() | Out-Null

# Add the new worksheet right after bingo_data; it becomes active/selected automatically
$ws = $wb.Worksheets.Add($null, $bingo)

$ws.Range("B1").Value = "names"
$ws.Range("C1").Value = "tas"
$ws.Range("D1").Value = "phases"
$ws.Range("E1").Value = "strat_buckets"
$ws.Range("F1").Value = "npvs"
$ws.Range("B2").Value = "Avniman"
$ws.Range("C2").Value = "CNS"
$ws.Range("D2").Value = "Phase 2"
$ws.Range("E2").Value = "Considered"
$ws.Range("F2").Value = 521.384600082794
$ws.Range("B3").Value = "Creficil"
$ws.Range("C3").Value = "Endocrine Disorder"
$ws.Range("D3").Value = "NDA"
$ws.Range("E3").Value = "Committed"
$ws.Range("F3").Value = 339.80008531156
$ws.Range("B4").Value = "Eaglogen"
$ws.Range("C4").Value = "Endocrine Disorder"
$ws.Range("D4").Value = "Phase 3"
$ws.Range("E4").Value = "Considered"
$ws.Range("F4").Value = 4736.33660931094
$ws.Range("B5").Value = "Estger"
$ws.Range("C5").Value = "Immunology"
$ws.Range("D5").Value = "Phase 1"
$ws.Range("E5").Value = "Considered"
$ws.Range("F5").Value = 902.40706268767
$ws.Range("B6").Value = "Holitorcitus"
$ws.Range("C6").Value = "Immunology"
$ws.Range("D6").Value = "Phase 2"
$ws.Range("E6").Value = "Considered"
$ws.Range("F6").Value = 2066.5810595901
$ws.Range("B7").Value = "Masogen"
$ws.Range("C7").Value = "Immunology"
$ws.Range("D7").Value = "Phase 1"
$ws.Range("E7").Value = "Potential"
$ws.Range("F7").Value = 376.33650758343
$ws.Range("B8").Value = "Matisem"
$ws.Range("C8").Value = "Endocrine Disorder"
$ws.Range("D8").Value = "Phase 1"
$ws.Range("E8").Value = "Considered"
$ws.Range("F8").Value = 1080.1330338124
$ws.Range("B9").Value = "Meprylol"
$ws.Range("C9").Value = "Dermatology"
$ws.Range("D9").Value = "Phase 1"
$ws.Range("E9").Value = "Considered"
$ws.Range("F9").Value = 268.81179113102
$ws.Range("B10").Value = "Mervisil"
$ws.Range("C10").Value = "Immunology"
$ws.Range("D10").Value = "NDA"
$ws.Range("E10").Value = "Committed"
$ws.Range("F10").Value = 358.41572150803
$ws.Range("B11").Value = "Metaphysis"
$ws.Range("C11").Value = "Immunology"
$ws.Range("D11").Value = "Phase 2"
$ws.Range("E11").Value = "Considered"
$ws.Range("F11").Value = 722.27320043144
$ws.Range("B12").Value = "Mrilipzor"
$ws.Range("C12").Value = "CNS"
$ws.Range("D12").Value = "Phase 3"
$ws.Range("E12").Value = "Considered"
$ws.Range("F12").Value = 1009.0816485826
$ws.Range("B13").Value = "Mritigen"
$ws.Range("C13").Value = "CNS"
$ws.Range("D13").Value = "Phase 1"
$ws.Range("E13").Value = "Considered"
$ws.Range("F13").Value = 688.86035319669
$ws.Range("B14").Value = "Nifilmox"
$ws.Range("C14").Value = "Ophthalmology"
$ws.Range("D14").Value = "Phase 2"
$ws.Range("E14").Value = "Considered"
$ws.Range("F14").Value = 1033.290529795
$ws.Range("B15").Value = "OpthTank"
$ws.Range("C15").Value = "Ophthalmology"
$ws.Range("D15").Value = "Phase 1"
$ws.Range("E15").Value = "Potential"
$ws.Range("F15").Value = 268.81179113102
$ws.Range("B16").Value = "Polgen"
$ws.Range("C16").Value = "Immunology"
$ws.Range("D16").Value = "Phase 3"
$ws.Range("E16").Value = "Committed"
$ws.Range("F16").Value = 632.37380423456
$ws.Range("B17").Value = "Prototase"
$ws.Range("C17").Value = "Endocrine Disorder"
$ws.Range("D17").Value = "Phase 2"
$ws.Range("E17").Value = "Considered"
$ws.Range("F17").Value = 910.03597482091
$ws.Range("B18").Value = "Refevel"
$ws.Range("C18").Value = "CNS"
$ws.Range("D18").Value = "Phase 2"
$ws.Range("E18").Value = "Considered"
$ws.Range("F18").Value = 267.69458691623
$ws.Range("B19").Value = "Reflitol"
$ws.Range("C19").Value = "CNS"
$ws.Range("D19").Value = "NDA"
$ws.Range("E19").Value = "Committed"
$ws.Range("F19").Value = 1047.067736859
$ws.Range("B20").Value = "Resdexel"
$ws.Range("C20").Value = "Endocrine Disorder"
$ws.Range("D20").Value = "Preclinical"
$ws.Range("E20").Value = "Potential"
$ws.Range("F20").Value = 528.271331410349
$ws.Range("B21").Value = "Rilopof"
$ws.Range("C21").Value = "CNS"
$ws.Range("D21").Value = "Phase 1"
$ws.Range("E21").Value = "Considered"
$ws.Range("F21").Value = 868.83527846219
$ws.Range("B22").Value = "Rydovanil"
$ws.Range("C22").Value = "CNS"
$ws.Range("D22").Value = "Phase 1"
$ws.Range("E22").Value = "Considered"
$ws.Range("F22").Value = 310.44710384782
$ws.Range("B23").Value = "Rytifil"
$ws.Range("C23").Value = "Ophthalmology"
$ws.Range("D23").Value = "NDA"
$ws.Range("E23").Value = "Committed"
$ws.Range("F23").Value = 361.39879135993
$ws.Range("B24").Value = "Tikofermin"
$ws.Range("C24").Value = "Dermatology"
$ws.Range("D24").Value = "Phase 3"
$ws.Range("E24").Value = "Committed"
$ws.Range("F24").Value = 1262.57665840502
$ws.Range("B25").Value = "Trivlexin"
$ws.Range("C25").Value = "Endocrine Disorder"
$ws.Range("D25").Value = "Phase 3"
$ws.Range("E25").Value = "Considered"
$ws.Range("F25").Value = 1024.4970446458
$ws.Range("B26").Value = "Varmenase"
$ws.Range("C26").Value = "Endocrine Disorder"
$ws.Range("D26").Value = "Phase 2"
$ws.Range("E26").Value = "Considered"
$ws.Range("F26").Value = 1694.5964688639
$ws.Range("B27").Value = "Virtiman"
$ws.Range("C27").Value = "Immunology"
$ws.Range("D27").Value = "Preclinical"
$ws.Range("E27").Value = "Considered"
$ws.Range("F27").Value = 350.47948177681
$ws.Range("B28").Value = "Vrexigen"
$ws.Range("C28").Value = "CNS"
$ws.Range("D28").Value = "Phase 1"
$ws.Range("E28").Value = "Potential"
$ws.Range("F28").Value = -8.81005425252839
$ws.Range("B29").Value = "Vrilimen"
$ws.Range("C29").Value = "Ophthalmology"
$ws.Range("D29").Value = "Phase 1"
$ws.Range("E29").Value = "Considered"
$ws.Range("F29").Value = 353.1369331391
$ws.Range("B30").Value = "Weglifil"
$ws.Range("C30").Value = "Immunology"
$ws.Range("D30").Value = "Phase 1"
$ws.Range("E30").Value = "Considered"
$ws.Range("F30").Value = 593.1908287328
$ws.Range("B31").Value = "Xumanase"
$ws.Range("C31").Value = "CNS"
$ws.Range("D31").Value = "Preclinical"
$ws.Range("E31").Value = "Considered"
$ws.Range("F31").Value = 1239.948635754
$ws.Range("B32").Value = "Xyfigil"
$ws.Range("C32").Value = "Endocrine Disorder"
$ws.Range("D32").Value = "Phase 1"
$ws.Range("E32").Value = "Potential"
$ws.Range("F32").Value = 3.3658485017045
$ws.Range("B33").Value = "Zerxil"
$ws.Range("C33").Value = "Endocrine Disorder"
$ws.Range("D33").Value = "Phase 3"
$ws.Range("E33").Value = "Considered"
$ws.Range("F33").Value = 172.21333019276

$ws.Range("H8").Select() | Out-Null
